$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "To-Test" flags for the two rows to YES
$ws.Range("C4").Value = "YES"
$ws.Range("C10").Value = "YES"

# Scroll the view back to the top and select C10 (as left by the editing user)
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C10").Select()
